# Append 32 new Lotofacil draw results (rows 424-455) to the LOTOFACIL sheet
# and update the active selection to match the post-edit view state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(3494,1,4,9,10,11,12,14,15,16,17,19,21,22,24,25),
    @(3495,2,3,6,8,9,11,12,15,16,18,19,20,22,23,25),
    @(3496,1,4,5,8,9,11,13,14,15,16,17,18,20,23,24),
    @(3497,1,2,3,4,6,10,11,12,13,14,15,18,19,21,25),
    @(3498,1,2,3,4,5,7,8,9,15,16,17,19,21,22,24),
    @(3499,1,5,7,9,10,11,12,13,14,17,18,19,23,24,25),
    @(3500,1,2,4,5,6,7,9,10,11,12,14,16,18,19,21),
    @(3501,1,3,7,9,11,12,13,15,16,18,19,20,21,22,24),
    @(3502,4,5,6,7,8,9,14,15,16,17,18,19,20,21,25),
    @(3503,1,4,6,7,8,10,11,12,14,15,16,18,21,24,25),
    @(3504,1,2,4,6,7,9,10,12,15,16,17,21,22,23,25),
    @(3505,1,2,3,4,6,7,8,9,11,14,16,20,21,23,25),
    @(3506,3,4,5,6,8,11,12,13,16,17,18,22,23,24,25),
    @(3507,2,3,4,7,9,10,11,12,17,18,19,21,22,24,25),
    @(3508,1,3,4,5,6,8,11,13,14,16,17,20,21,22,24),
    @(3509,1,3,4,7,9,10,12,14,15,16,19,21,22,24,25),
    @(3510,1,2,3,4,5,6,7,8,13,14,15,17,21,22,23),
    @(3511,2,3,5,6,7,8,9,12,13,14,15,19,20,23,24),
    @(3512,1,2,3,4,5,6,9,12,13,14,15,17,19,20,22),
    @(3513,3,4,5,6,7,11,13,14,15,16,18,20,22,23,24),
    @(3514,1,3,4,5,6,8,10,12,13,16,17,20,21,23,24),
    @(3515,1,2,3,4,5,6,7,9,14,15,16,17,21,23,24),
    @(3516,2,3,6,11,12,13,14,15,17,18,19,20,21,23,25),
    @(3517,1,2,3,7,8,11,12,14,17,18,19,20,21,22,24),
    @(3518,2,5,6,7,11,12,13,14,16,17,19,21,22,24,25),
    @(3519,4,5,6,7,8,13,14,16,17,18,20,21,22,23,25),
    @(3520,1,3,7,8,10,12,13,15,16,18,19,20,21,22,23),
    @(3521,6,7,8,9,10,11,12,13,14,15,19,20,21,22,23),
    @(3522,1,2,4,5,7,8,10,11,12,14,15,17,20,22,24),
    @(3523,5,6,8,9,10,11,12,14,15,16,17,19,23,24,25),
    @(3524,2,3,4,5,7,9,12,13,14,15,16,17,18,22,24),
    @(3525,1,2,3,7,10,11,14,15,16,17,18,20,21,23,24)
)

$startRow = 424
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# Scroll the frozen pane toward the bottom of the new data and select the
# newly added block, matching the workbook state after entering the draws.
$excel.ActiveWindow.ScrollRow = 418
$null = $ws.Range("B435:P455").Select()
